# Update gh-pages to output generated at 456a3b4
# Applies updated "want to go" counts and a refreshed cover image URL
# to the "展览" (Exhibition) and "全部类型" (All types) sheets, plus a
# single count update on the "演出" (Performance) sheet.

$wb = $excel.ActiveWorkbook

$newImageUrl = "//i2.hdslb.com/bfs/openplatform/202410/bZ8StmGX1729756956769.jpeg"

# --- Sheet "展览" (Exhibition) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 31
$wsExhibit.Range("I3").Value = $newImageUrl
$wsExhibit.Range("F6").Value = 5280
$wsExhibit.Range("F8").Value = 93
$wsExhibit.Range("F9").Value = 102
$wsExhibit.Range("F10").Value = 361
$wsExhibit.Range("F11").Value = 14

# --- Sheet "演出" (Performance) ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F6").Value = 128

# --- Sheet "全部类型" (All types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 31
$wsAll.Range("I3").Value = $newImageUrl
$wsAll.Range("F9").Value = 5280
$wsAll.Range("F11").Value = 93
$wsAll.Range("F12").Value = 102
$wsAll.Range("F14").Value = 361
$wsAll.Range("F15").Value = 14
$wsAll.Range("F17").Value = 128
